# Generate Report for Handback
# Adds a new handback record (9fcba326-d7a1-4eab-9d69-bf6d1fade614) to the
# Overview, zh-cn and de-de sheets - one new row appended to each table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)
$ovRow = $ovTable.ListRows.Add()
$ovR = $ovRow.Range().Row()

$ov.Cells.Item($ovR, 1).Value = "9fcba326-d7a1-4eab-9d69-bf6d1fade614.md"
$ov.Cells.Item($ovR, 3).Value = ".md"
$ov.Cells.Item($ovR, 5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item($ovR, 6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item($ovR, 7).Value = "2016-10-21 03:39:01"
$ov.Cells.Item($ovR, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Cells.Item($ovR, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae5c00e2aea01e600087215fa145b6349fd01208/e2e/9fcba326-d7a1-4eab-9d69-bf6d1fade614.md", "", "", "e2e\9fcba326-d7a1-4eab-9d69-bf6d1fade614.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhRow = $zhTable.ListRows.Add()
$zhR = $zhRow.Range().Row()

$zh.Cells.Item($zhR, 2).Value = ".md"
$zh.Cells.Item($zhR, 3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item($zhR, 4).Value = "e2e"
$zh.Cells.Item($zhR, 5).Value = "ht"
$zh.Cells.Item($zhR, 6).Value = "True"
$zh.Cells.Item($zhR, 7).Value = "9fcba326-d7a1-4eab-9d69-bf6d1fade614.08f1fb9741fe8abc9badef2f976bc71375143004.zh-cn.xlf"
$zh.Cells.Item($zhR, 8).Value = "2016-10-21 03:38:50"
$zh.Cells.Item($zhR, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item($zhR, 10).Value = "9fcba326-d7a1-4eab-9d69-bf6d1fade614.08f1fb9741fe8abc9badef2f976bc71375143004.zh-cn.xlf"
$zh.Cells.Item($zhR, 11).Value = "2016-10-21 03:39:32"
$zh.Cells.Item($zhR, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Cells.Item($zhR, 12).Value = ""
$zh.Cells.Item($zhR, 13).Value = "True"
$zh.Cells.Item($zhR, 14).Value = ""
$zh.Cells.Item($zhR, 15).Value = "False"
$zh.Cells.Item($zhR, 16).Value = ""

$zh.Hyperlinks.Add($zh.Cells.Item($zhR, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c62bb2f7688b36c2d929b1e031336852c0d929a7/e2e/9fcba326-d7a1-4eab-9d69-bf6d1fade614.md", "", "", "9fcba326-d7a1-4eab-9d69-bf6d1fade614.md")
$zh.Hyperlinks.Add($zh.Cells.Item($zhR, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c62bb2f7688b36c2d929b1e031336852c0d929a7/e2e/9fcba326-d7a1-4eab-9d69-bf6d1fade614.md", "", "", "9fcba326-d7a1-4eab-9d69-bf6d1fade614.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deRow = $deTable.ListRows.Add()
$deR = $deRow.Range().Row()

$de.Cells.Item($deR, 2).Value = ".md"
$de.Cells.Item($deR, 3).Value = "Handed back: in sync with en-US"
$de.Cells.Item($deR, 4).Value = "e2e"
$de.Cells.Item($deR, 5).Value = "ht"
$de.Cells.Item($deR, 6).Value = "True"
$de.Cells.Item($deR, 7).Value = "9fcba326-d7a1-4eab-9d69-bf6d1fade614.08f1fb9741fe8abc9badef2f976bc71375143004.de-de.xlf"
$de.Cells.Item($deR, 8).Value = "2016-10-21 03:39:01"
$de.Cells.Item($deR, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item($deR, 10).Value = "9fcba326-d7a1-4eab-9d69-bf6d1fade614.08f1fb9741fe8abc9badef2f976bc71375143004.de-de.xlf"
$de.Cells.Item($deR, 11).Value = "2016-10-21 03:39:50"
$de.Cells.Item($deR, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Cells.Item($deR, 12).Value = ""
$de.Cells.Item($deR, 13).Value = "True"
$de.Cells.Item($deR, 14).Value = ""
$de.Cells.Item($deR, 15).Value = "False"
$de.Cells.Item($deR, 16).Value = ""

$de.Hyperlinks.Add($de.Cells.Item($deR, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d0eaa2cd726684b36cead3914f8fb6645197d52e/e2e/9fcba326-d7a1-4eab-9d69-bf6d1fade614.md", "", "", "9fcba326-d7a1-4eab-9d69-bf6d1fade614.md")
$de.Hyperlinks.Add($de.Cells.Item($deR, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d0eaa2cd726684b36cead3914f8fb6645197d52e/e2e/9fcba326-d7a1-4eab-9d69-bf6d1fade614.md", "", "", "9fcba326-d7a1-4eab-9d69-bf6d1fade614.md")

Write-Output "Handback report row added for 9fcba326-d7a1-4eab-9d69-bf6d1fade614"
